$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Hours" column (B2:B5) stores its numbers as text (shared strings),
# not numeric values. Update each cell to the new text value, using a
# leading apostrophe so Excel keeps storing it as text instead of
# auto-converting the numeric-looking string into a real number.
$ws.Range("B2").Value = "'518.0640000000001"
$ws.Range("B3").Value = "'946.193"
$ws.Range("B4").Value = "'325.302"
$ws.Range("B5").Value = "'1352.583"

# Drop the quote-prefix formatting that Excel applies when text is forced
# via a leading apostrophe, so the cells keep their original (unstyled)
# appearance.
$ws.Range("B2:B5").ClearFormats()
